# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 08:30"

# --- Swap "Islas Malvinas" / "Montserrat" rows (216 / 217) ---
# Row 216 becomes Montserrat, row 217 becomes Islas Malvinas
$ws.Range("A216").Value = "Montserrat"
$ws.Range("A217").Value = "Islas Malvinas"

# Their Casos activos (D) / Muertes (H) values swap along with the name swap
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0

# --- Row 28: Ucrania ---
$ws.Range("B28").Value = 293641
$ws.Range("C28").Value = 6410
$ws.Range("D28").Value = 124113
$ws.Range("E28").Value = 164011
$ws.Range("G28").Value = 109
$ws.Range("H28").Value = 5517

# --- Row 68: Kirguistan ---
$ws.Range("B68").Value = 51490
$ws.Range("C68").Value = 470
$ws.Range("D68").Value = 45509
$ws.Range("E68").Value = 4873
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 1108

# --- Row 83: El Salvador ---
$ws.Range("E83").Value = 3806
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 917
